# Update "want to go" counts (column F) on sheet "展览" (index 1)
# and sheet "全部类型" (index 4) to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value = 75
$wsExpo.Range("F3").Value = 131
$wsExpo.Range("F4").Value = 2084
$wsExpo.Range("F5").Value = 367
$wsExpo.Range("F6").Value = 632
$wsExpo.Range("F7").Value = 102
$wsExpo.Range("F8").Value = 2074
$wsExpo.Range("F9").Value = 10699
$wsExpo.Range("F12").Value = 286
$wsExpo.Range("F15").Value = 7555
$wsExpo.Range("F18").Value = 260
$wsExpo.Range("F20").Value = 3337

# --- Sheet 4: 全部类型 ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 75
$wsAll.Range("F3").Value = 131
$wsAll.Range("F4").Value = 2084
$wsAll.Range("F5").Value = 367
$wsAll.Range("F6").Value = 632
$wsAll.Range("F8").Value = 102
$wsAll.Range("F9").Value = 2074
$wsAll.Range("F12").Value = 10699
$wsAll.Range("F15").Value = 286
$wsAll.Range("F18").Value = 7555
$wsAll.Range("F21").Value = 260
$wsAll.Range("F23").Value = 3337
